{"js": "// Add two new paragraphs right after the paragraph that ends with\n// \"I manually labeled speaker gender\", pushing the existing trailing\n// empty paragraph further down (unchanged).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst anchorText = \"I manually labeled speaker gender\";\nconst anchor = paragraphs.items.find((p) => p.text === anchorText);\nif (!anchor) {\n  throw new Error(\"Could not find anchor paragraph: \" + anchorText);\n}\n\n// Insert in order, each time right after the previous inserted paragraph.\nconst p1 = anchor.insertParagraph(\n  \"Speakers 1-10 are mandarin speakers, 11-20 are English speakers\",\n  \"After\"\n);\np1.insertParagraph(\n  \"As a native speaker of both languages, it really sounds like PRC Mandarin and US English to me\",\n  \"After\"\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$anchorText = \"I manually labeled speaker gender\"\n\n$anchor = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq $anchorText) {\n        $anchor = $p\n        break\n    }\n}\n\nif ($anchor -eq $null) {\n    Write-Output \"Anchor paragraph not found: $anchorText\"\n} else {\n    $cr = [char]13\n    $newText = $cr + \"Speakers 1-10 are mandarin speakers, 11-20 are English speakers\" + $cr + \"As a native speaker of both languages, it really sounds like PRC Mandarin and US English to me\"\n    $anchor.Range.InsertAfter($newText)\n    Write-Output \"Inserted paragraphs after anchor.\"\n}\n"}
